# Apply updated cryptos list values (prices / 1h volume%) per the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.801.17"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.622.29"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'596.30"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "'149.78"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("D11").Value = "'5.61"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "'27.54"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "3.094.93"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "63.683.00"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "'0.0000149"
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("D17").Value = "2.607.48"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("E18").Value = "  +6.45%  "
$ws.Range("D19").Value = "'4.60"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "'348.53"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").Value = "'6.89"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'5.72"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("D24").Value = "'66.18"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("E25").Value = "  +13.19%  "
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").Value = "'9.17"
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'8.11"
$ws.Range("E28").Value = "  +3.50%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.164"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "'541.24"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("D33").Value = "0.0₃0843"
$ws.Range("E33").Value = "  +4.70%  "
$ws.Range("D34").Value = "'1.74"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("D35").Value = "'5.20"
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("D36").Value = "'168.65"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  +2.32%  "
$ws.Range("D40").Value = "'19.36"
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'169.78"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").Value = "'39.79"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "'3.91"
$ws.Range("E44").Value = "  +3.84%  "
$ws.Range("D45").Value = "'0.0598"
$ws.Range("E45").Value = "  +3.42%  "
$ws.Range("D46").Value = "'21.33"
$ws.Range("E46").Value = "  -5.50%  "
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").Value = "'1.97"
$ws.Range("E49").Value = "  +10.33%  "
$ws.Range("D50").Value = "'0.0967"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "'19.16"
$ws.Range("E51").Value = "  +1.47%  "
